$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 7 de Mayo de 2020 a las 10:34"

# Filipinas overtook Dinamarca in total cases, so the two rows swap order
# and Filipinas' stats are refreshed with the new totals.
$ws.Range("A42").Value = "Filipinas"
$ws.Range("B42").Value = 10343
$ws.Range("C42").Value = 339
$ws.Range("D42").Value = 1618
$ws.Range("E42").Value = 8040
$ws.Range("F42").Value = 31
$ws.Range("G42").Value = 27
$ws.Range("H42").Value = 685

$ws.Range("A43").Value = "Dinamarca"
$ws.Range("B43").Value = 10083
$ws.Range("C43").Value = 145
$ws.Range("D43").Value = 7493
$ws.Range("E43").Value = 2084
$ws.Range("F43").Value = 46
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 506

# Estonia (row 81) daily stat refresh
$ws.Range("B81").Value = 1720
$ws.Range("C81").Value = 7
$ws.Range("D81").Value = 273
$ws.Range("E81").Value = 1391
$ws.Range("G81").Value = 1
$ws.Range("H81").Value = 56

# Sri Lanka (row 102) daily stat refresh
$ws.Range("D102").Value = 232
$ws.Range("E102").Value = 556
